$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update theater names to drop trailing descriptors
$ws.Range("A1").Value = "1. AMC Empire"
$ws.Range("A5").Value = "5. Cinemark Tinseltown"
$ws.Range("A6").Value = "6. Cinemark Playa Vista"

# Move the active selection to A12, as recorded in the saved file
$ws.Range("A12").Select()
